# TradingModel - 2021/11/15 data update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated data rows (StockIdx, StockId, PositionSize) for rows 2-10
$data = @(
    @(1,  2436, 60),
    @(3,  3035, 32),
    @(4,  3122, 90),
    @(5,  3141, 27),
    @(8,  3588, 35),
    @(10, 6104, 36),
    @(11, 6138, 30),
    @(13, 6271, 20),
    @(14, 6411, 26)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}
